$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string must be forced to
# remain Text (matching the source data, which always stores these columns
# as strings) - otherwise Excel auto-converts "0.380" -> 0.38, etc.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '60.711.54'
$ws.Range('D3').Value = '2.640.69'
$ws.Range('E3').Value = '  +1.27%  '
$ws.Range('E4').Value = '  +0.01%  '
Set-TextValue $ws.Range('D5') '576.59'
$ws.Range('E5').Value = '  -0.41%  '
Set-TextValue $ws.Range('D6') '143.79'
$ws.Range('E6').Value = '  +0.04%  '
Set-TextValue $ws.Range('D7') '0.998'
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  -0.47%  '
$ws.Range('E9').Value = '  +0.63%  '
$ws.Range('E10').Value = '  +0.11%  '
Set-TextValue $ws.Range('D11') '0.380'
$ws.Range('E11').Value = '  +2.46%  '
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('D13').Value = '3.110.11'
$ws.Range('E13').Value = '  +1.10%  '
$ws.Range('E14').Value = '  +11.61%  '
$ws.Range('D15').Value = '60.693.99'
$ws.Range('E15').Value = '  -0.23%  '
$ws.Range('E16').Value = '  -0.31%  '
$ws.Range('D17').Value = '2.655.92'
$ws.Range('E17').Value = '  +1.19%  '
Set-TextValue $ws.Range('D18') '11.54'
$ws.Range('E18').Value = '  +2.22%  '
Set-TextValue $ws.Range('D19') '4.71'
$ws.Range('E19').Value = '  +0.91%  '
Set-TextValue $ws.Range('D20') '349.78'
$ws.Range('E20').Value = '  -0.17%  '
$ws.Range('E21').Value = '  -1.08%  '
$ws.Range('E22').Value = '  +0.11%  '
Set-TextValue $ws.Range('D23') '0.527'
$ws.Range('E23').Value = '  +1.77%  '
Set-TextValue $ws.Range('D24') '63.79'
$ws.Range('E24').Value = '  +0.77%  '
Set-TextValue $ws.Range('D25') '0.997'
$ws.Range('E25').Value = '  +0.02%  '
Set-TextValue $ws.Range('D26') '0.161'
$ws.Range('E26').Value = '  +0.26%  '
Set-TextValue $ws.Range('D27') '8.18'
$ws.Range('E28').Value = '  +9.45%  '
$ws.Range('D29').Value = '0.0₃0803'
$ws.Range('E29').Value = '  +0.26%  '
$ws.Range('E30').Value = '  +6.60%  '
$ws.Range('E31').Value = '  +0.09%  '
Set-TextValue $ws.Range('D32') '163.23'
$ws.Range('E32').Value = '  +0.44%  '
$ws.Range('E33').Value = '  +1.52%  '
Set-TextValue $ws.Range('D34') '4.61'
$ws.Range('E35').Value = '  +3.63%  '
Set-TextValue $ws.Range('D36') '1.32'
$ws.Range('E36').Value = '  +6.77%  '
$ws.Range('B37').Value = 'Stacks'
$ws.Range('C37').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D37') '1.66'
$ws.Range('E37').Value = '  +2.18%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range('D38') '338.84'
$ws.Range('E38').Value = '  +9.66%  '
Set-TextValue $ws.Range('D40') '0.904'
$ws.Range('E40').Value = '  +6.65%  '
Set-TextValue $ws.Range('D41') '38.39'
$ws.Range('E41').Value = '  +1.10%  '
$ws.Range('E42').Value = '  +2.87%  '
Set-TextValue $ws.Range('D43') '0.623'
$ws.Range('E43').Value = '  +2.22%  '
Set-TextValue $ws.Range('D44') '20.26'
$ws.Range('E44').Value = '  +1.46%  '
$ws.Range('E45').Value = '  +2.74%  '
$ws.Range('E46').Value = '  +2.18%  '
Set-TextValue $ws.Range('D47') '132.80'
$ws.Range('E47').Value = '  -0.87%  '
$ws.Range('E48').Value = '  +0.91%  '
Set-TextValue $ws.Range('D49') '20.55'
$ws.Range('E49').Value = '  +0.05%  '
Set-TextValue $ws.Range('D50') '0.998'
$ws.Range('E50').Value = '  +0.30%  '
$ws.Range('D51').Value = '2.085.61'
$ws.Range('E51').Value = '  +2.32%  '
